# "adding validation for authentication"
#
# The "Tanggal Lahir" (date of birth) column (J) holds dates as plain
# text. They were entered ambiguously as MM/DD/YYYY, which clashes with
# day values > 12 (e.g. 01/13/2015) and is easy to mis-parse downstream.
# Normalize every value in that column to an unambiguous YYYY/MM/DD text
# form so validation logic can rely on a single, sortable, ISO-like
# layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp           = [Microsoft.Office.Interop.Excel.XlDirection]::xlUp
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$dateCol = 10  # column J - "Tanggal Lahir"

$lastRow = $ws.Cells.Item($ws.Rows.Count, $dateCol).End($xlUp).Row
if ($lastRow -lt 2) { $lastRow = 45 }

$dataRange = $ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol))

# Switch the whole column to text format first so the reassigned
# "YYYY/MM/DD" strings are kept as literal text instead of being
# auto-recognized and silently converted into date serial numbers.
$dataRange.NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    $val = $cell.Value2
    if ($val -ne $null -and $val -match '^(\d{2})/(\d{2})/(\d{4})$') {
        $mm = $Matches[1]
        $dd = $Matches[2]
        $yyyy = $Matches[3]
        $cell.Value2 = "$yyyy/$mm/$dd"
    }
}

# Restore the column's original look (right-aligned, General-formatted,
# same as every other untouched cell) by pasting the formatting back in
# from the never-modified header cell above the data.
$template = $ws.Cells.Item(1, $dateCol)
$template.Copy() | Out-Null
$dataRange.PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
